# "Add Lair of the Magi blueprint"
# Adds 13 new enemy rows (230-242) to the "Enemies" sheet, and makes
# "Enemies" the active/selected sheet (previously "Armour" was selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enemies")

# Data for the new rows, in sheet (top-to-bottom) order.
# Columns: Name(A/B/C), D..P numeric stats (all 0 for these new blueprint
# monsters), Q/R notes ("{}" placeholder), S/T, U (fleeable), V (deathHp-ish
# counter), W (music cue), X/Y (living/flags).
$rows = @(
    @{ Row=230; Name="Henchman";         V=0; W="None";             X=0; Y=1 },
    @{ Row=231; Name="Arcane Sorcerer";   V=0; W="None";             X=0; Y=1 },
    @{ Row=232; Name="Death Warrior";     V=0; W="None";             X=1; Y=0 },
    @{ Row=233; Name="Death Magus";       V=0; W="None";             X=2; Y=0 },
    @{ Row=234; Name="Fire Zealot";       V=0; W="None";             X=3; Y=1 },
    @{ Row=235; Name="Abomination";       V=0; W="None";             X=4; Y=1 },
    @{ Row=236; Name="Conjurer";          V=0; W="None";             X=5; Y=1 },
    @{ Row=237; Name="Conjured Demon";    V=0; W="None";             X=6; Y=1 },
    @{ Row=238; Name="Undead Magus";      V=0; W="None";             X=7; Y=0 },
    @{ Row=239; Name="Magus Wraith";      V=0; W="None";             X=8; Y=0 },
    @{ Row=240; Name="Time Wizard";       V=0; W="None";             X=9; Y=1 },
    @{ Row=241; Name="Touin DePenk";      V=2; W="Important Battle"; X=1; Y=1 },
    @{ Row=242; Name="Wonnen Daztinque";  V=1; W="Important Battle"; X=1; Y=1 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Name   # A - identifier
    $ws.Cells.Item($row, 2).Value = $r.Name   # B - name
    $ws.Cells.Item($row, 3).Value = $r.Name   # C - image key

    $ws.Cells.Item($row, 4).Value = 0         # D - level
    $ws.Cells.Item($row, 5).Value = 0         # E - hp
    $ws.Cells.Item($row, 6).Value = 0         # F - maxHp
    $ws.Cells.Item($row, 7).Value = 0         # G - damage
    $ws.Cells.Item($row, 8).Value = 0         # H - accuracy
    $ws.Cells.Item($row, 9).Value = 0         # I - cRate
    $ws.Cells.Item($row, 10).Value = 0        # J - cDamage
    $ws.Cells.Item($row, 11).Value = 0        # K - bRate
    $ws.Cells.Item($row, 12).Value = 0        # L - defence
    $ws.Cells.Item($row, 13).Value = 0        # M - earthReduction
    $ws.Cells.Item($row, 14).Value = 0        # N - waterReduction
    $ws.Cells.Item($row, 15).Value = 0        # O - fireReduction
    $ws.Cells.Item($row, 16).Value = 0        # P - physicalReduction

    $ws.Cells.Item($row, 17).Value = "{}"     # Q - SKILLS
    $ws.Cells.Item($row, 18).Value = "{}"     # R - ITEMS

    $ws.Cells.Item($row, 19).Value = 0        # S - XP
    $ws.Cells.Item($row, 20).Value = 0        # T - EUROS
    $ws.Cells.Item($row, 21).Value = 1        # U - FLEEABLE
    $ws.Cells.Item($row, 22).Value = $r.V     # V - DEATH_HP
    $ws.Cells.Item($row, 23).Value = $r.W     # W - MUSIC
    $ws.Cells.Item($row, 24).Value = $r.X     # X - UNIQUE
    $ws.Cells.Item($row, 25).Value = $r.Y     # Y - LIVING
}

# Make "Enemies" the active sheet (was "Armour"); restore the frozen-pane
# layout (header row + 4 id columns) and leave the final selection on the
# first newly-added row, matching the edited workbook's view state.
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("E2").Select()
$win.FreezePanes = $true
$ws.Range("A230").Select()
